$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '23.171.50'
$ws.Range("E2").Value = '  +0.30%  '
$ws.Range("D3").Value = "'" + '1.600.21'
$ws.Range("E3").Value = '  -0.04%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("E5").Value = '  +0.04%  '
$ws.Range("D6").Value = "'" + '303.13'
$ws.Range("E6").Value = '  +0.49%  '
$ws.Range("D7").Value = "'" + '0.3782'
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = "'" + '52.15'
$ws.Range("E8").Value = '  +4.63%  '
$ws.Range("D9").Value = "'" + '0.3612'
$ws.Range("E9").Value = '  -1.12%  '
$ws.Range("D10").Value = "'" + '1.264'
$ws.Range("E10").Value = '  -0.50%  '
$ws.Range("E11").Value = '  +0.11%  '
$ws.Range("D12").Value = "'" + '0.08117'
$ws.Range("E12").Value = '  -0.52%  '
$ws.Range("D13").Value = "'" + '22.64'
$ws.Range("E13").Value = '  -2.16%  '
$ws.Range("D14").Value = "'" + '6.575'
$ws.Range("E14").Value = '  -0.41%  '
$ws.Range("D15").Value = "'" + '7.394'
$ws.Range("E15").Value = '  +0.18%  '
$ws.Range("E16").Value = '  -1.16%  '
$ws.Range("D17").Value = "'" + '1.601.35'
$ws.Range("E17").Value = '  -0.06%  '
$ws.Range("D18").Value = "'" + '93.99'
$ws.Range("E18").Value = '  +2.63%  '
$ws.Range("D19").Value = "'" + '0.06899'
$ws.Range("E19").Value = '  +0.63%  '
$ws.Range("D20").Value = "'" + '18.05'
$ws.Range("E20").Value = '  -1.67%  '
$ws.Range("D21").Value = "'" + '6.540'
$ws.Range("E21").Value = '  -0.46%  '
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").Value = "'" + '12.95'
$ws.Range("E23").Value = '  -0.54%  '
$ws.Range("D24").Value = "'" + '23.167.05'
$ws.Range("E24").Value = '  +0.31%  '
$ws.Range("D25").Value = "'" + '2.379'
$ws.Range("E25").Value = '  +1.60%  '
$ws.Range("D26").Value = "'" + '2.998'
$ws.Range("E26").Value = '  +9.84%  '
$ws.Range("D27").Value = "'" + '21.24'
$ws.Range("E27").Value = '  +0.44%  '
$ws.Range("D28").Value = "'" + '149.48'
$ws.Range("E28").Value = '  -0.66%  '
$ws.Range("D29").Value = "'" + '5.250'
$ws.Range("E29").Value = '  -0.55%  '
$ws.Range("D30").Value = "'" + '133.64'
$ws.Range("E30").Value = '  +0.96%  '
$ws.Range("D31").Value = "'" + '2.385'
$ws.Range("E31").Value = '  -0.69%  '
$ws.Range("D32").Value = "'" + '6.799'
$ws.Range("E32").Value = '  -0.62%  '
$ws.Range("D33").Value = "'" + '1.780.35'
$ws.Range("E33").Value = '  -0.50%  '
$ws.Range("D34").Value = "'" + '0.9669'
$ws.Range("E34").Value = '  +0.49%  '
$ws.Range("D35").Value = "'" + '0.07486'
$ws.Range("E35").Value = '  -2.36%  '
$ws.Range("D36").Value = "'" + '10.29'
$ws.Range("E36").Value = '  +1.92%  '
$ws.Range("E37").Value = '  -0.98%  '
$ws.Range("D38").Value = "'" + '0.2508'
$ws.Range("E38").Value = '  -1.88%  '
$ws.Range("D39").Value = "'" + '0.08803'
$ws.Range("E39").Value = '  -1.11%  '
$ws.Range("D40").Value = "'" + '6.084'
$ws.Range("E40").Value = '  -3.31%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = "'" + '1.360'
$ws.Range("E41").Value = '  -0.71%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = "'" + '0.7095'
$ws.Range("E42").Value = '  -0.14%  '
$ws.Range("E43").Value = '  -1.73%  '
$ws.Range("D44").Value = "'" + '15.51'
$ws.Range("E44").Value = '  +0.32%  '
$ws.Range("E45").Value = '  -1.64%  '
$ws.Range("D46").Value = "'" + '2.310'
$ws.Range("E46").Value = '  -0.24%  '
$ws.Range("D47").Value = "'" + '4.011'
$ws.Range("E47").Value = '  +0.61%  '
$ws.Range("D48").Value = "'" + '131.99'
$ws.Range("E48").Value = '  -0.22%  '
$ws.Range("D49").Value = "'" + '0.07957'
$ws.Range("E49").Value = '  +0.27%  '
$ws.Range("D50").Value = "'" + '1.201'
$ws.Range("E50").Value = '  -0.54%  '
$ws.Range("E51").Value = '  +1.83%  '
